# Update gh-pages to output generated at 456a3b4
# Applies the scraped data refresh to the "展览" (Worksheets index 1) and
# "全部类型" (Worksheets index 4) sheets.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# ---- "展览" sheet updates ----
$wsExhibit.Range("F3").Value  = 1525
$wsExhibit.Range("F4").Value  = 904
$wsExhibit.Range("F11").Value = 2477
$wsExhibit.Range("G14").Value = 218
$wsExhibit.Range("F16").Value = 621
$wsExhibit.Range("F17").Value = 795
$wsExhibit.Range("F24").Value = 5098
$wsExhibit.Range("F26").Value = 558
$wsExhibit.Range("F27").Value = 84
$wsExhibit.Range("F33").Value = 1044
$wsExhibit.Range("F36").Value = 57
$wsExhibit.Range("F39").Value = 1078
$wsExhibit.Range("F44").Value = 55

# ---- "全部类型" sheet updates ----
$wsAll.Range("F5").Value  = 1526
$wsAll.Range("F6").Value  = 904
$wsAll.Range("F17").Value = 2477
$wsAll.Range("G20").Value = 218
$wsAll.Range("F22").Value = 621
$wsAll.Range("F24").Value = 795
$wsAll.Range("F29").Value = 5098
$wsAll.Range("F31").Value = 558
$wsAll.Range("F32").Value = 84
$wsAll.Range("F38").Value = 1044
$wsAll.Range("F40").Value = 57
$wsAll.Range("F42").Value = 1078
$wsAll.Range("F46").Value = 55
